# Apply "working on update user settings" edits to the 用户设置 (Sheet3) worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Row 12: the update-user-info API/method becomes update-user-profile
$ws.Range("D12").Value = "api - updateUserProfile"
$ws.Range("G12").Value = "method - update_user_profile"

# Row 15: "账号管理Tab" becomes "密码Tab", and a new controller label appears at F15
$ws.Range("C15").Value = "密码Tab"
$ws.Range("F15").Value = "controller - Users"

# Row 16: updateEmail is replaced by updatePassword (+ its method label)
$ws.Range("D16").Value = "api - updatePassword"
$ws.Range("G16").Value = "method - update_password"

# Row 17 (old "api - updatePassword") is no longer used
$ws.Range("D17").ClearContents()

# Row 18 stays "api - updatePhone" (already correct, no change needed)

# Row 19 (new): updateEmail moves down here
$ws.Range("D19").Value = "api - updateEmail"

# Make 用户设置 (Sheet3) the active sheet/tab, with G21 selected
$ws.Activate()
$ws.Range("G21").Select()
